$wb = $excel.ActiveWorkbook

$mainSheet = $wb.Worksheets.Item("Main")
$setupSheet = $wb.Worksheets.Item("__SETUP__")

# ---------------------------------------------------------------------
# 1. Duplicate the "Main" sheet, placing the copy right after "Main", and
#    rename it to "Abort". This gives us the new abort-sequence tab with
#    the same layout/validation/conditional-formatting as the main
#    sequence tab.
# ---------------------------------------------------------------------
$mainSheet.Copy($null, $mainSheet)
$abortSheet = $wb.Worksheets.Item("Main (2)")
$abortSheet.Name = "Abort"

# ---------------------------------------------------------------------
# 2. Re-point the sheet-scoped named ranges so the Abort sheet gets its
#    own TestSeqTimes/TestSeqValves/TestSeqActions that feed from its own
#    columns, while keeping the workbook/Main-scoped names intact (just
#    recreated so they land after the Abort-scoped ones, matching the
#    authored ordering).
# ---------------------------------------------------------------------
$wb.Names.Item("TestSeqActions").Delete()
$abortSheet.Names.Add("TestSeqActions", "=Abort!`$D`$4:`$D`$1048576")
$wb.Names.Add("TestSeqActions", "=Main!`$D`$4:`$D`$1048576")

$wb.Names.Item("TestSeqTimes").Delete()
$abortSheet.Names.Add("TestSeqTimes", "=Abort!`$B`$4:`$B`$1048576")
$wb.Names.Add("TestSeqTimes", "=Main!`$B`$4:`$B`$1048576")

$wb.Names.Item("TestSeqValves").Delete()
$abortSheet.Names.Add("TestSeqValves", "=Abort!`$C`$4:`$C`$1048576")
$mainSheet.Names.Add("TestSeqValves", "=Main!`$C`$4:`$C`$1048576")

# ---------------------------------------------------------------------
# 3. Populate the Abort sheet with its own trigger/sensor-port sequence
#    data (decoded from the abort test-sequence file) instead of the
#    values it inherited from Main.
# ---------------------------------------------------------------------
$abortSheet.Range("B4").Value = 9
$abortSheet.Range("B5").Value = 9
$abortSheet.Range("B6").Value = 10
$abortSheet.Range("B7").Value = 10
$abortSheet.Range("B8").Value = 2
$abortSheet.Range("B9").Value = 2

$abortSheet.Range("F4").Value = 12
$abortSheet.Range("G4").Value = 12
$abortSheet.Range("F5").Value = 13
$abortSheet.Range("G5").Value = 14

$abortSheet.Range("H4").Value = "FPD-1"
$abortSheet.Range("H5").Value = "FPD-2"

# ---------------------------------------------------------------------
# 4. The Main sheet's trigger sensor-port mapping is also updated (the
#    sensor port numbering moved from FPD-3 placeholders to the correct
#    FPD-1/FPD-2 ports).
# ---------------------------------------------------------------------
$mainSheet.Range("H4").Value = "FPD-1"
$mainSheet.Range("H5").Value = "FPD-2"

# ---------------------------------------------------------------------
# 5. Restore each sheet's selection/active-cell to match the authored
#    session, and make Abort the active tab.
# ---------------------------------------------------------------------
$mainSheet.Range("E17").Select()
$abortSheet.Range("H8").Select()
$setupSheet.Range("F5").Select()

$abortSheet.Activate()
